# Update top marker genes per cluster sheet with pts-filtered results
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("MeV.1.4.1")
$ws.Range("A5").Value = "Clec14a"
$ws.Range("B5").Value = 9.15392780303955
$ws.Range("C5").Value = 4.631550312042236
$ws.Range("D5").Value = [double]"5.680524071573974E-18"
$ws.Range("E5").Value = 0.3476702508960574
$ws.Range("A6").Value = "Ptprb"
$ws.Range("B6").Value = 24.30417060852051
$ws.Range("C6").Value = 4.537421226501465
$ws.Range("D6").Value = [double]"9.254136655594219E-127"
$ws.Range("E6").Value = 0.942652329749104

$ws = $wb.Worksheets.Item("MeV.1.4.7")
$ws.Range("A3").Value = "Igfbp3"
$ws.Range("B3").Value = 5.603856086730957
$ws.Range("C3").Value = 2.349907636642456
$ws.Range("D3").Value = [double]"3.71320565257575E-06"
$ws.Range("E3").Value = 0.3024390243902439
$ws.Range("A4").Value = "Pla1a"
$ws.Range("B4").Value = 6.086055755615234
$ws.Range("C4").Value = 2.349749326705933
$ws.Range("D4").Value = [double]"2.573184392021096E-07"
$ws.Range("E4").Value = 0.3121951219512195
$ws.Range("A5").Value = "Ranbp3l"
$ws.Range("B5").Value = 12.89632034301758
$ws.Range("C5").Value = 2.10603666305542
$ws.Range("D5").Value = [double]"2.466879160583189E-34"
$ws.Range("E5").Value = 0.8926829268292683
$ws.Range("A6").Value = "Slc7a11"
$ws.Range("B6").Value = 13.03956604003906
$ws.Range("C6").Value = 2.0783531665802
$ws.Range("D6").Value = [double]"7.425888279816912E-35"
$ws.Range("E6").Value = 0.8634146341463415

$ws = $wb.Worksheets.Item("MeV.1.4.8")
$ws.Range("A2").Value = "Ftl1"
$ws.Range("B2").Value = 7.98644495010376
$ws.Range("C2").Value = 3.640774965286255
$ws.Range("D2").Value = [double]"7.257124971693554E-13"
$ws.Range("E2").Value = 0.3823529411764706
$ws.Range("A3").Value = "Rpl18a"
$ws.Range("B3").Value = 8.461723327636719
$ws.Range("C3").Value = 3.634722232818604
$ws.Range("D3").Value = [double]"2.202588431015998E-14"
$ws.Range("E3").Value = 0.4019607843137255
$ws.Range("A4").Value = "Rpl19"
$ws.Range("B4").Value = 7.77732515335083
$ws.Range("C4").Value = 3.607499599456787
$ws.Range("D4").Value = [double]"3.365697566386239E-12"
$ws.Range("E4").Value = 0.3725490196078431
$ws.Range("A5").Value = "Rpl41"
$ws.Range("B5").Value = 10.50026988983154
$ws.Range("C5").Value = 3.580030679702759
$ws.Range("D5").Value = [double]"1.800275661234721E-22"
$ws.Range("E5").Value = 0.5098039215686274
$ws.Range("A6").Value = "Rps23"
$ws.Range("B6").Value = 8.355746269226074
$ws.Range("C6").Value = 3.565969705581665
$ws.Range("D6").Value = [double]"4.686185945897644E-14"
$ws.Range("E6").Value = 0.4019607843137255

$ws = $wb.Worksheets.Item("MeV.2.1")
$ws.Range("A4").Value = "Grem2"
$ws.Range("B4").Value = 7.478861331939697
$ws.Range("C4").Value = 2.485360860824585
$ws.Range("D4").Value = [double]"3.482079853079695E-11"
$ws.Range("E4").Value = 0.35
$ws.Range("A5").Value = "Aox3"
$ws.Range("B5").Value = 7.725160121917725
$ws.Range("C5").Value = 2.30659556388855
$ws.Range("D5").Value = [double]"5.559300436662962E-12"
$ws.Range("E5").Value = 0.3714285714285714
$ws.Range("A6").Value = "Tmem132c"
$ws.Range("B6").Value = 6.671573162078857
$ws.Range("C6").Value = 2.137815475463867
$ws.Range("D6").Value = [double]"8.396090083369488E-09"
$ws.Range("E6").Value = 0.3357142857142857

$ws = $wb.Worksheets.Item("MeV.2.8")
$ws.Range("A2").Value = "Tmem132e"
$ws.Range("B2").Value = 7.111780166625977
$ws.Range("C2").Value = 3.38364577293396
$ws.Range("D2").Value = [double]"1.961631211988044E-10"
$ws.Range("E2").Value = 0.3383838383838384
$ws.Range("A3").Value = "Igf1"
$ws.Range("B3").Value = 17.9183349609375
$ws.Range("C3").Value = 3.092049598693848
$ws.Range("D3").Value = [double]"8.865119719549452E-68"
$ws.Range("E3").Value = 0.9393939393939394
$ws.Range("A4").Value = "Ltbp1"
$ws.Range("B4").Value = 16.49800491333008
$ws.Range("C4").Value = 3.040301084518433
$ws.Range("D4").Value = [double]"1.981494228641063E-57"
$ws.Range("E4").Value = 0.8484848484848485
$ws.Range("A5").Value = "Grem2"
$ws.Range("B5").Value = 10.15969562530518
$ws.Range("C5").Value = 2.973413467407227
$ws.Range("D5").Value = [double]"1.306358055235668E-21"
$ws.Range("E5").Value = 0.51010101010101
$ws.Range("A6").Value = "Svep1"
$ws.Range("B6").Value = 13.78193855285645
$ws.Range("C6").Value = 2.91600775718689
$ws.Range("D6").Value = [double]"5.263196384214766E-40"
$ws.Range("E6").Value = 0.7121212121212122

$ws = $wb.Worksheets.Item("MeV.3.17")
$ws.Range("A6").Value = "Tbx18"
$ws.Range("B6").Value = 5.165728092193604
$ws.Range("C6").Value = 4.254248142242432
$ws.Range("D6").Value = [double]"6.764717505333729E-05"
$ws.Range("E6").Value = 0.3444444444444444

$ws = $wb.Worksheets.Item("MeV.4.1")
$ws.Range("A4").Value = "Gm12002"
$ws.Range("B4").Value = 7.802682399749756
$ws.Range("C4").Value = 6.187079429626465
$ws.Range("D4").Value = [double]"8.279112683654652E-13"
$ws.Range("E4").Value = 0.3439153439153439
$ws.Range("A5").Value = "Notch3"
$ws.Range("B5").Value = 19.91693878173828
$ws.Range("C5").Value = 5.95043420791626
$ws.Range("D5").Value = [double]"1.010888609432314E-84"
$ws.Range("E5").Value = 0.8783068783068783
$ws.Range("A6").Value = "Tbx3os1"
$ws.Range("B6").Value = 15.91995525360107
$ws.Range("C6").Value = 5.90712833404541
$ws.Range("D6").Value = [double]"4.376945966447012E-54"
$ws.Range("E6").Value = 0.7037037037037037

$ws = $wb.Worksheets.Item("MeV.4.12")
$ws.Range("A2").Value = "Myoc"
$ws.Range("B2").Value = 5.7233567237854
$ws.Range("C2").Value = 5.037930011749268
$ws.Range("D2").Value = [double]"1.254535432636743E-06"
$ws.Range("E2").Value = 0.3303571428571428
$ws.Range("A3").Value = "Gm973"
$ws.Range("B3").Value = 6.428671360015869
$ws.Range("C3").Value = 4.575037956237793
$ws.Range("D3").Value = [double]"2.118475017243618E-08"
$ws.Range("E3").Value = 0.3839285714285715
$ws.Range("A4").Value = "Prps2"
$ws.Range("B4").Value = 6.665574550628662
$ws.Range("C4").Value = 4.530786037445068
$ws.Range("D4").Value = [double]"4.669635598598761E-09"
$ws.Range("E4").Value = 0.3928571428571428
$ws.Range("A5").Value = "Slc47a1"
$ws.Range("B5").Value = 16.18927574157715
$ws.Range("C5").Value = 4.450317859649658
$ws.Range("D5").Value = [double]"1.254710437143614E-54"
$ws.Range("E5").Value = 0.9910714285714286
$ws.Range("A6").Value = "Tbx15"
$ws.Range("B6").Value = 10.17861366271973
$ws.Range("C6").Value = 4.307960510253906
$ws.Range("D6").Value = [double]"1.475351188655028E-21"
$ws.Range("E6").Value = 0.6071428571428571

$ws = $wb.Worksheets.Item("MeV.4.21")
$ws.Range("A2").Value = "Sema3g"
$ws.Range("B2").Value = 5.77971363067627
$ws.Range("C2").Value = 7.014009952545166
$ws.Range("D2").Value = [double]"1.109204878479648E-06"
$ws.Range("E2").Value = 0.3636363636363636
$ws.Range("A3").Value = "Bmx"
$ws.Range("B3").Value = 10.23504638671875
$ws.Range("C3").Value = 6.608002662658691
$ws.Range("D3").Value = [double]"2.062170293983623E-21"
$ws.Range("E3").Value = 0.6477272727272727
$ws.Range("A4").Value = "Nos1"
$ws.Range("B4").Value = 5.670414447784424
$ws.Range("C4").Value = 5.827220439910889
$ws.Range("D4").Value = [double]"2.011758290660017E-06"
$ws.Range("E4").Value = 0.3636363636363636
$ws.Range("A5").Value = "Prdm16"
$ws.Range("B5").Value = 10.99948310852051
$ws.Range("C5").Value = 5.73930025100708
$ws.Range("D5").Value = [double]"8.925358263575392E-25"
$ws.Range("E5").Value = 0.7045454545454546
$ws.Range("A6").Value = "Lama3"
$ws.Range("B6").Value = 7.815247058868408
$ws.Range("C6").Value = 4.904033184051514
$ws.Range("D6").Value = [double]"2.204871793142434E-12"
$ws.Range("E6").Value = 0.5113636363636364

$ws = $wb.Worksheets.Item("MeV.1.4.11")
$ws.Range("A6").Value = "Ppp1r1a"
$ws.Range("B6").Value = 5.995941162109375
$ws.Range("C6").Value = 2.513557195663452
$ws.Range("D6").Value = [double]"7.550876496636515E-07"
$ws.Range("E6").Value = 0.3409090909090909

$ws = $wb.Worksheets.Item("MeV.4.30")
$ws.Range("A3").Value = "Angptl1"
$ws.Range("B3").Value = 4.304681777954102
$ws.Range("C3").Value = 5.369338035583496
$ws.Range("D3").Value = 0.002859872853422683
$ws.Range("E3").Value = 0.3518518518518519
$ws.Range("A4").Value = "Nox4"
$ws.Range("B4").Value = 7.003499031066895
$ws.Range("C4").Value = 4.927291393280029
$ws.Range("D4").Value = [double]"3.06934761341701E-09"
$ws.Range("E4").Value = 0.5925925925925926
$ws.Range("A5").Value = "Col14a1"
$ws.Range("B5").Value = 5.971248626708984
$ws.Range("C5").Value = 4.874123096466064
$ws.Range("D5").Value = [double]"1.366950806119628E-06"
$ws.Range("E5").Value = 0.5
$ws.Range("A6").Value = "Gfpt2"
$ws.Range("B6").Value = 3.939765214920044
$ws.Range("C6").Value = 4.800660610198975
$ws.Range("D6").Value = 0.01065446690907652
$ws.Range("E6").Value = 0.3333333333333333

$ws = $wb.Worksheets.Item("MeV.1.4.15")
$ws.Range("A3").Value = "Slco1c1"
$ws.Range("B3").Value = 4.474020481109619
$ws.Range("C3").Value = 2.42079496383667
$ws.Range("D3").Value = 0.0005093369476376944
$ws.Range("E3").Value = 0.3217391304347826
$ws.Range("A4").Value = "Spock2"
$ws.Range("B4").Value = 3.735859632492065
$ws.Range("C4").Value = 2.281996488571167
$ws.Range("D4").Value = 0.007636826084559388
$ws.Range("E4").Value = 0.3130434782608696
$ws.Range("A5").Value = "Adgrl4"
$ws.Range("B5").Value = 6.849327564239502
$ws.Range("C5").Value = 2.203131437301636
$ws.Range("D5").Value = [double]"2.349714339750008E-09"
$ws.Range("E5").Value = 0.5304347826086957
$ws.Range("A6").Value = "Tmtc2"
$ws.Range("B6").Value = 8.723294258117676
$ws.Range("C6").Value = 2.179532766342163
$ws.Range("D6").Value = [double]"2.97259115242356E-15"
$ws.Range("E6").Value = 0.7217391304347827

$ws = $wb.Worksheets.Item("MeV.1.4.2")
$ws.Range("A3").Value = "Cntnap5b"
$ws.Range("B3").Value = 3.528460264205933
$ws.Range("C3").Value = 1.252711772918701
$ws.Range("D3").Value = 0.02151797131806221
$ws.Range("E3").Value = 0.356
$ws.Range("A4").Value = "Gria3"
$ws.Range("B4").Value = 3.638855218887329
$ws.Range("C4").Value = 1.027507901191711
$ws.Range("D4").Value = 0.01621328852190851
$ws.Range("E4").Value = 0.412
$ws.Range("A5").Value = "Kirrel3"
$ws.Range("B5").Value = 6.181618690490723
$ws.Range("C5").Value = 1.027142882347107
$ws.Range("D5").Value = [double]"3.489781664528542E-07"
$ws.Range("E5").Value = 0.616
$ws.Range("A6").Value = "Lsamp"
$ws.Range("B6").Value = 10.45835208892822
$ws.Range("C6").Value = 1.014389038085938
$ws.Range("D6").Value = [double]"2.80427248342351E-21"
$ws.Range("E6").Value = 0.948

$ws = $wb.Worksheets.Item("MeV.1.4.5")
$ws.Range("A5").Value = "Tfrc"
$ws.Range("B5").Value = 8.758108139038086
$ws.Range("C5").Value = 3.983325242996216
$ws.Range("D5").Value = [double]"2.545998775480264E-16"
$ws.Range("E5").Value = 0.3699186991869919
$ws.Range("A6").Value = "Nostrin"
$ws.Range("B6").Value = 9.850622177124023
$ws.Range("C6").Value = 3.98274040222168
$ws.Range("D6").Value = [double]"1.095224244319312E-20"
$ws.Range("E6").Value = 0.4105691056910569

